$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 4.919465525594682
$ws.Cells.Item(2, 4).Value = 4.409400734803915
$ws.Cells.Item(2, 5).Value = 16.4951395720656
$ws.Cells.Item(2, 6).Value = 26.6947121690564
$ws.Cells.Item(2, 7).Value = 34.01207126782955
$ws.Cells.Item(2, 8).Value = 14.77919364037958
$ws.Cells.Item(2, 11).Value = 15.41918535586897
$ws.Cells.Item(2, 14).Value = 16.57566166230366
$ws.Cells.Item(3, 3).Value = 4.749625503947096
$ws.Cells.Item(3, 4).Value = 4.435624843057894
$ws.Cells.Item(3, 5).Value = 15.55565519777343
$ws.Cells.Item(3, 6).Value = 26.45266724202762
$ws.Cells.Item(3, 7).Value = 33.46534592053941
$ws.Cells.Item(3, 8).Value = 14.78013410565199
$ws.Cells.Item(3, 11).Value = 14.79378306681102
$ws.Cells.Item(3, 14).Value = 16.64864034559091
$ws.Cells.Item(4, 3).Value = 4.643851622340542
$ws.Cells.Item(4, 4).Value = 4.452294301199195
$ws.Cells.Item(4, 5).Value = 14.95454649616157
$ws.Cells.Item(4, 6).Value = 26.31402890675283
$ws.Cells.Item(4, 7).Value = 33.14204103004788
$ws.Cells.Item(4, 8).Value = 14.78521115709388
$ws.Cells.Item(4, 11).Value = 14.39942447985864
$ws.Cells.Item(4, 14).Value = 16.69539907217806
$ws.Cells.Item(5, 3).Value = 4.600449782367975
$ws.Cells.Item(5, 4).Value = 4.459230989532921
$ws.Cells.Item(5, 5).Value = 14.70375955856703
$ws.Cells.Item(5, 6).Value = 26.26009492591661
$ws.Cells.Item(5, 7).Value = 33.01359274275764
$ws.Cells.Item(5, 8).Value = 14.78840669533324
$ws.Cells.Item(5, 11).Value = 14.23637088457637
$ws.Cells.Item(5, 14).Value = 16.71494594496888
$ws.Cells.Item(6, 3).Value = 4.593227297196554
$ws.Cells.Item(6, 4).Value = 4.460391532135322
$ws.Cells.Item(6, 5).Value = 14.6617734021746
$ws.Cells.Item(6, 6).Value = 26.2512954107495
$ws.Cells.Item(6, 7).Value = 32.99246863778009
$ws.Cells.Item(6, 8).Value = 14.78900519815208
$ws.Cells.Item(6, 11).Value = 14.20916210841163
$ws.Cells.Item(6, 14).Value = 16.71822148411764
$ws.Cells.Item(7, 3).Value = 4.643267391014173
$ws.Cells.Item(7, 4).Value = 4.452387268532384
$ws.Cells.Item(7, 5).Value = 14.95118750546409
$ws.Cells.Item(7, 6).Value = 26.31329109450434
$ws.Cells.Item(7, 7).Value = 33.14029512665609
$ws.Cells.Item(7, 8).Value = 14.78524969885075
$ws.Cells.Item(7, 11).Value = 14.39723463514613
$ws.Cells.Item(7, 14).Value = 16.6956606920917
$ws.Cells.Item(8, 3).Value = 4.861258513574546
$ws.Cells.Item(8, 4).Value = 4.418325566711061
$ws.Cells.Item(8, 5).Value = 16.17638735136423
$ws.Cells.Item(8, 6).Value = 26.60921793352317
$ws.Cells.Item(8, 7).Value = 33.82110406046517
$ws.Cells.Item(8, 8).Value = 14.77858133007537
$ws.Cells.Item(8, 11).Value = 15.20584113224446
$ws.Cells.Item(8, 14).Value = 16.60042135100358
$ws.Cells.Item(9, 3).Value = 5.273835501643875
$ws.Cells.Item(9, 4).Value = 4.355989912697232
$ws.Cells.Item(9, 5).Value = 18.46581625186574
$ws.Cells.Item(9, 6).Value = 27.26606377911448
$ws.Cells.Item(9, 7).Value = 35.24586130336765
$ws.Cells.Item(9, 8).Value = 14.80138925311246
$ws.Cells.Item(9, 11).Value = 16.6995262267186
$ws.Cells.Item(9, 14).Value = 16.42902964644935
$ws.Cells.Item(10, 3).Value = 5.564202007861616
$ws.Cells.Item(10, 4).Value = 4.312845476745442
$ws.Cells.Item(10, 5).Value = 20.12041306727535
$ws.Cells.Item(10, 6).Value = 27.79138892426623
$ws.Cells.Item(10, 7).Value = 36.33544699304316
$ws.Cells.Item(10, 8).Value = 14.84021964656524
$ws.Cells.Item(10, 11).Value = 17.72961336871365
$ws.Cells.Item(10, 14).Value = 16.31234360073391
$ws.Cells.Item(11, 3).Value = 5.692876379534704
$ws.Cells.Item(11, 4).Value = 4.293780069879173
$ws.Cells.Item(11, 5).Value = 20.83139714743529
$ws.Cells.Item(11, 6).Value = 28.03872663214829
$ws.Cells.Item(11, 7).Value = 36.83787202225562
$ws.Cells.Item(11, 8).Value = 14.86270175730886
$ws.Cells.Item(11, 11).Value = 18.18171910370141
$ws.Cells.Item(11, 14).Value = 16.26123641431803
$ws.Cells.Item(12, 3).Value = 5.741064629336696
$ws.Cells.Item(12, 4).Value = 4.286640105443743
$ws.Cells.Item(12, 5).Value = 21.09467149190681
$ws.Cells.Item(12, 6).Value = 28.13350722493543
$ws.Cells.Item(12, 7).Value = 37.0288898893991
$ws.Cells.Item(12, 8).Value = 14.87190839855944
$ws.Cells.Item(12, 11).Value = 18.35042018614659
$ws.Cells.Item(12, 14).Value = 16.2421650966214
$ws.Cells.Item(13, 3).Value = 5.730711041176981
$ws.Cells.Item(13, 4).Value = 4.28817429636613
$ws.Cells.Item(13, 5).Value = 21.03823513896629
$ws.Cells.Item(13, 6).Value = 28.11304610901524
$ws.Cells.Item(13, 7).Value = 36.98772029024003
$ws.Cells.Item(13, 8).Value = 14.86989475984153
$ws.Cells.Item(13, 11).Value = 18.31420053135857
$ws.Cells.Item(13, 14).Value = 16.24625994144792
$ws.Cells.Item(14, 3).Value = 5.696851886115235
$ws.Cells.Item(14, 4).Value = 4.293191069384052
$ws.Cells.Item(14, 5).Value = 20.8531761225618
$ws.Cells.Item(14, 6).Value = 28.04650234602771
$ws.Cells.Item(14, 7).Value = 36.85357313890292
$ws.Cells.Item(14, 8).Value = 14.86344530274763
$ws.Cells.Item(14, 11).Value = 18.19564906085347
$ws.Cells.Item(14, 14).Value = 16.2596617678343
$ws.Cells.Item(15, 3).Value = 5.676040822143462
$ws.Cells.Item(15, 4).Value = 4.296274335797457
$ws.Cells.Item(15, 5).Value = 20.73904705934887
$ws.Cells.Item(15, 6).Value = 28.00588559638099
$ws.Cells.Item(15, 7).Value = 36.7714969300639
$ws.Cells.Item(15, 8).Value = 14.85958508737514
$ws.Cells.Item(15, 11).Value = 18.12270352473643
$ws.Cells.Item(15, 14).Value = 16.26790742207084
$ws.Cells.Item(16, 3).Value = 5.555719984817951
$ws.Cells.Item(16, 4).Value = 4.314102652705289
$ws.Cells.Item(16, 5).Value = 20.07311259758604
$ws.Cells.Item(16, 6).Value = 27.77538604598504
$ws.Cells.Item(16, 7).Value = 36.30273022601885
$ws.Cells.Item(16, 8).Value = 14.83884738732474
$ws.Cells.Item(16, 11).Value = 17.69972383001114
$ws.Cells.Item(16, 14).Value = 16.31572312581245
$ws.Cells.Item(17, 3).Value = 5.480997293100894
$ws.Cells.Item(17, 4).Value = 4.325182770446168
$ws.Cells.Item(17, 5).Value = 19.65392899078725
$ws.Cells.Item(17, 6).Value = 27.636065707013
$ws.Cells.Item(17, 7).Value = 36.01673760126592
$ws.Cells.Item(17, 8).Value = 14.82736000850906
$ws.Cells.Item(17, 11).Value = 17.43591775493421
$ws.Cells.Item(17, 14).Value = 16.34556065902135
$ws.Cells.Item(18, 3).Value = 5.437699797381657
$ws.Cells.Item(18, 4).Value = 4.331608645910558
$ws.Cells.Item(18, 5).Value = 19.40889846068353
$ws.Cells.Item(18, 6).Value = 27.55672484446622
$ws.Cells.Item(18, 7).Value = 35.85289167066279
$ws.Cells.Item(18, 8).Value = 14.82120630075453
$ws.Cells.Item(18, 11).Value = 17.28263826344462
$ws.Cells.Item(18, 14).Value = 16.36290832748139
$ws.Cells.Item(19, 3).Value = 5.422986744408851
$ws.Cells.Item(19, 4).Value = 4.333793450623699
$ws.Cells.Item(19, 5).Value = 19.32525971580651
$ws.Cells.Item(19, 6).Value = 27.53000001454659
$ws.Cells.Item(19, 7).Value = 35.79753434938188
$ws.Cells.Item(19, 8).Value = 14.8192006296449
$ws.Cells.Item(19, 11).Value = 17.23047937190302
$ws.Cells.Item(19, 14).Value = 16.36881394030243
$ws.Cells.Item(20, 3).Value = 5.488985021670985
$ws.Cells.Item(20, 4).Value = 4.323997806928531
$ws.Cells.Item(20, 5).Value = 19.69895798681445
$ws.Cells.Item(20, 6).Value = 27.65081517141625
$ws.Cells.Item(20, 7).Value = 36.0471163422133
$ws.Cells.Item(20, 8).Value = 14.82853591575612
$ws.Cells.Item(20, 11).Value = 17.4641612941865
$ws.Cells.Item(20, 14).Value = 16.34236517621491
$ws.Cells.Item(21, 3).Value = 5.706812086845672
$ws.Cells.Item(21, 4).Value = 4.291715366813571
$ws.Cells.Item(21, 5).Value = 20.90769390094317
$ws.Cells.Item(21, 6).Value = 28.06601816631945
$ws.Cells.Item(21, 7).Value = 36.89295644335663
$ws.Cells.Item(21, 8).Value = 14.86532085451644
$ws.Cells.Item(21, 11).Value = 18.2305393425721
$ws.Cells.Item(21, 14).Value = 16.25571769252829
$ws.Cells.Item(22, 3).Value = 5.846021169102611
$ws.Cells.Item(22, 4).Value = 4.271080955644329
$ws.Cells.Item(22, 5).Value = 21.66295070851986
$ws.Cells.Item(22, 6).Value = 28.34386045815825
$ws.Cells.Item(22, 7).Value = 37.45009913619729
$ws.Cells.Item(22, 8).Value = 14.89340167582091
$ws.Cells.Item(22, 11).Value = 18.71678290574766
$ws.Cells.Item(22, 14).Value = 16.20073056641998
$ws.Cells.Item(23, 3).Value = 5.772025299276865
$ws.Cells.Item(23, 4).Value = 4.282051805642975
$ws.Cells.Item(23, 5).Value = 21.26302096500816
$ws.Cells.Item(23, 6).Value = 28.1950057818731
$ws.Cells.Item(23, 7).Value = 37.15241370380666
$ws.Cells.Item(23, 8).Value = 14.87804490346018
$ws.Cells.Item(23, 11).Value = 18.45864257105117
$ws.Cells.Item(23, 14).Value = 16.22992863630553
$ws.Cells.Item(24, 3).Value = 5.485374820085267
$ws.Cells.Item(24, 4).Value = 4.324533355153623
$ws.Cells.Item(24, 5).Value = 19.6786129526246
$ws.Cells.Item(24, 6).Value = 27.64414457667284
$ws.Cells.Item(24, 7).Value = 36.03338030612592
$ws.Cells.Item(24, 8).Value = 14.82800288479669
$ws.Cells.Item(24, 11).Value = 17.45139740679242
$ws.Cells.Item(24, 14).Value = 16.34380925059348
$ws.Cells.Item(25, 3).Value = 5.164220994474449
$ws.Cells.Item(25, 4).Value = 4.372382661170157
$ws.Cells.Item(25, 5).Value = 17.81927425216374
$ws.Cells.Item(25, 6).Value = 27.0805614134352
$ws.Cells.Item(25, 7).Value = 34.85203856617607
$ws.Cells.Item(25, 8).Value = 14.79135361952661
$ws.Cells.Item(25, 11).Value = 16.30655152107677
$ws.Cells.Item(25, 14).Value = 16.47376402841673
